$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new columns I (I0) and J (IF) per row
$values = @{
    2  = @(8, 9)
    3  = @(5, 7)
    4  = @(3, 6)
    5  = @(8, 9)
    6  = @(8, 8)
    7  = @(2, 5)
    8  = @(5, 7)
    9  = @(11, 11)
    10 = @(7, 7)
    11 = @(5, 6)
    12 = @(8, 9)
    13 = @(8, 9)
    14 = @(5, 5)
    15 = @(5, 6)
    16 = @(5, 6)
    17 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
